$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Tyler Anderson - update Last.Updated date and Injury.Details text
$ws.Range("C2").Value = "July 27 2017"
$ws.Range("E2").Value = "Anderson has been transferred to the 60-day disabled list due to arthroscopic left knee surgery and is expected to be sidelined until the end of August."

# Row 3 (Chad Bettis) remains unchanged

# Row 4: replace Tyler Chatwood entry with Ian Desmond entry
$ws.Range("A4").Value = "Ian Desmond"
$ws.Range("B4").Value = "desmoia01"
$ws.Range("C4").Value = "July 28 2017"
$ws.Range("D4").Value = "Calf"
$ws.Range("E4").Value = "Desmond has been placed on the 10-day disabled due to a right calf strain. A timetable for his recovery has yet to be established."

# Row 5: replace Gerardo Parra entry with Jake McGee entry
$ws.Range("A5").Value = "Jake McGee"
$ws.Range("B5").Value = "mcgeeja01"
$ws.Range("C5").Value = "July 30 2017"
$ws.Range("D5").Value = "Back"
$ws.Range("E5").Value = "McGee has been placed on the 10-day disabled list with a back injury and it is unclear how much time he is expected to miss."

# Update the active selection to A15
$ws.Range("A15").Select()
